$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10: clear the milestone/completed marks (previously II / X)
$ws.Range("E10:F10").ClearContents()

# Row 22: mark as Milestone III completed
$ws.Range("E22").Value = "III"
$ws.Range("F22").Value = "X"

# Row 67: mark as Milestone III completed
$ws.Range("E67").Value = "III"
$ws.Range("F67").Value = "X"

# Row 68: mark as Milestone III completed
$ws.Range("E68").Value = "III"
$ws.Range("F68").Value = "X"

# Update the selected/active cell and scroll position for the sheet view
$ws.Activate()
$ws.Range("F39").Select()
